$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Version bump 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date bump
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher now populated
$meta.Range("B9").Value = "Alvearie Team"

# Old duplicate "Contact" rows (10 & 11) become "Jurisdiction" and the
# (already-present) "Description" row simply moves up once row 11 is removed.
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"
$meta.Rows.Item(11).Delete()

# --- Elements sheet updates ---
$elem = $wb.Worksheets.Item("Elements")

# Root extension row: Short/Definition replaced with the real age-group text
$elem.Range("K2").Value = "Age Group"
$elem.Range("L2").Value = "Standard code for age groupings"
